$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.986.67"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").Value = "3.437.78"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'408.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").Value = "'128.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("E7").Value = "  +5.56%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.736"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.86%  "
$ws.Range("D10").Value = "'0.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +18.83%  "
$ws.Range("D11").Value = "'42.59"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.141"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.964.34"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'21.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.85%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.0000212"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +65.20%  "
$ws.Range("D16").Value = "'8.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.62%  "
$ws.Range("D17").Value = "3.369.58"
$ws.Range("E17").Value = "  -1.15%  "
$ws.Range("D18").Value = "'12.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +12.53%  "
$ws.Range("E19").Value = "  +5.53%  "
$ws.Range("D20").Value = "61.923.02"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").Value = "'402.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +26.68%  "
$ws.Range("D22").Value = "'89.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.44%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'13.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.67%  "
$ws.Range("D25").Value = "'3.22"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").Value = "'32.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.64%  "
$ws.Range("D27").Value = "'8.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.64%  "
$ws.Range("D28").Value = "'4.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "'7.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D31").Value = "'0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'43.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.04%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'11.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.27%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'0.0504"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("D37").Value = "'53.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("D39").Value = "'3.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'2.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.04%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "'0.132"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.34%  "
$ws.Range("E42").Value = "  +6.88%  "
$ws.Range("D43").Value = "'142.18"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "'4.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  +8.66%  "
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'21.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").Value = "2.121.25"
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'0.130"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.36%  "
$ws.Range("D51").Value = "'0.0378"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.08%  "
